$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 21.01.2022 18:00"

# Row 7 (MOL Olomoucka): swap price values, update delta and timestamp
$ws.Range("B7").Value = 36.7
$ws.Range("C7").Value = 36.9
# Leading apostrophe forces the delta to stay text (e.g. "-0.2") instead of
# being auto-converted to a number, matching the original "+0.2" text cell.
$ws.Range("D7").Value = "'-0.2"
$ws.Range("E7").Value = "2022-01-21 18:01:51"
